$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update shared strings / header text via Characters (preserves run formatting position) ----
# A8: "Volume 30   Number  9" -> "Volume 30   Number  10"
$c = $ws.Range("A8")
$c.Characters(21, 1).Text = "10"

# C9: "Report Covering the Week  2/27/2023  Through  3/5/2023"
#     -> "Report Covering the Week  3/6/2023  Through  3/12/2023"
$c = $ws.Range("C9")
$c.Characters(47, 8).Text = "3/12/2023"
$c.Characters(27, 9).Text = "3/6/2023"

# ---- Update crime statistics table (rows 15-30) ----

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = -50
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -80
$ws.Range("L15").Value = -50
$ws.Range("M15").Value = -50
$ws.Range("N15").Value = -66.666666666666

# Row 16
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -62.5
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -27.777777777777
$ws.Range("I16").Value = 34
$ws.Range("J16").Value = 39
$ws.Range("K16").Value = -12.820512820512
$ws.Range("L16").Value = 47.826086956521
$ws.Range("M16").Value = 3.030303030303
$ws.Range("N16").Value = -82.383419689119

# Row 17
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 21.428571428571
$ws.Range("I17").Value = 43
$ws.Range("J17").Value = 44
$ws.Range("K17").Value = -2.272727272727
$ws.Range("L17").Value = 22.857142857142
$ws.Range("M17").Value = 152.941176470588
$ws.Range("N17").Value = -10.416666666666

# Row 18
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 38
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = -15.555555555555
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -45.714285714285
$ws.Range("N18").Value = -88.790560471976

# Row 19
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 23
$ws.Range("E19").Value = -21.739130434782
$ws.Range("F19").Value = 76
$ws.Range("G19").Value = 76
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 182
$ws.Range("J19").Value = 180
$ws.Range("K19").Value = 1.111111111111
$ws.Range("L19").Value = 52.941176470588
$ws.Range("M19").Value = -30.798479087452
$ws.Range("N19").Value = -64.243614931237

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 0
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 10
$ws.Range("J20").Value = 10
$ws.Range("L20").Value = 233.333333333333
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = -96.015936254980

# Row 21
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -35
$ws.Range("F21").Value = 125
$ws.Range("G21").Value = 130
$ws.Range("H21").Value = -3.846153846153
$ws.Range("I21").Value = 308
$ws.Range("J21").Value = 323
$ws.Range("K21").Value = -4.643962848297
$ws.Range("L21").Value = 40
$ws.Range("M21").Value = -21.025641025641
$ws.Range("N21").Value = -77.100371747211

# Row 22
$ws.Range("C22").Value = 2
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 3
$ws.Range("E22").Value = -33.333333333333
$ws.Range("F22").Value = 10
$ws.Range("G22").Value = 8
$ws.Range("H22").Value = 25
$ws.Range("I22").Value = 19
$ws.Range("J22").Value = 23
$ws.Range("K22").Value = -17.391304347826
$ws.Range("L22").Value = 137.5
$ws.Range("M22").Value = 46.153846153846

# Row 24
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 43
$ws.Range("E24").Value = -11.627906976744
$ws.Range("F24").Value = 135
$ws.Range("G24").Value = 190
$ws.Range("H24").Value = -28.947368421052
$ws.Range("I24").Value = 377
$ws.Range("J24").Value = 435
$ws.Range("K24").Value = -13.333333333333
$ws.Range("L24").Value = 8.333333333333
$ws.Range("M24").Value = 12.202380952381

# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 9.090909090909
$ws.Range("F25").Value = 35
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = -22.222222222222
$ws.Range("I25").Value = 83
$ws.Range("J25").Value = 89
$ws.Range("K25").Value = -6.741573033707
$ws.Range("L25").Value = 45.614035087719
$ws.Range("M25").Value = 15.277777777777

# Row 26
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 7
$ws.Range("K26").Value = -42.857142857142
$ws.Range("L26").Value = 100

# Row 27
$ws.Range("D27").Value = 3
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 66.666666666666
$ws.Range("E27").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F27").Value = 11
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 37.5
$ws.Range("I27").Value = 20
$ws.Range("J27").Value = 17
$ws.Range("K27").Value = 17.647058823529
$ws.Range("L27").Value = 53.846153846153

# Row 30
$ws.Range("D30").Value = 2
$ws.Range("G30").Value = 6
$ws.Range("J30").Value = 7
$ws.Range("K30").Value = -85.714285714285
